$d = $word.ActiveDocument

# Use Find (search-only) + Range.Text assignment rather than Find.Execute's
# built-in Replace parameter -- the latter runs the replacement text through
# AutoCorrect/AutoFormat (turning straight apostrophes into curly ones),
# while a plain Range.Text assignment inserts the literal text untouched.
function Replace-Text($find, $replace) {
    $rng = $d.Content
    $found = $rng.Find.Execute($find)
    if (-not $found) {
        Write-Output "NOT FOUND: $find"
        return
    }
    $rng.Text = $replace
}

# --- Title ---
Replace-Text "Quirks of Quantum Realm" "History: Unveiling the Tapestry of Our Past"

# --- Author name: merge 3 runs into one ---
$pName = $d.Paragraphs.Item(2).Range
$pName.MoveEnd(1, -1) | Out-Null
$pName.Text = "Marcus Montgomery"

# --- Author email: rebuild as marcusmontgomery999@gmail.com ---
Replace-Text "alexander" "marcusmontgomery999@gmail"
Replace-Text ".maxwell@emailhub." "."
Replace-Text "net" "com"

# --- Body paragraph (quantum -> history) ---
Replace-Text "In the realm of quantum mechanics, paradoxical phenomena reign supreme, challenging our conventional understanding of the universe" "In every corner of the world, a story unfolds, a story that spans centuries and civilizations, a story of triumph and adversity woven together in a tapestry of history"

Replace-Text " Quantum entanglement and superposition beckon us to question the very nature of reality" " As we delve into its depths, we embark upon a journey that stretches far beyond textbooks, enriching our understanding of the present and illuminating the path to a better future"

Replace-Text " Entanglement epitomizes the spooky action at a distance, where particles separated by vast cosmic gulfs remain eerily correlated, defying Einstein's cherished idea of locality" " It is within this tapestry that we find the echoes of our ancestors' voices, whispered through artifacts and texts, engravings and ruins"

Replace-Text " And superposition presents a vision of boundless possibility, where particles can exist in multiple states simultaneously, transcending the limits of classical intuition. Unveiling the secrets of the quantum realm holds the promise of transformative technologies, from superfast quantum computers to cryptographic systems impervious to eavesdropping" " Each thread of this tapestry tells a tale of bravery and tragedy, of wisdom and folly, captivating our imaginations and inviting us to scrutinize humanity's complex and ever-changing nature"

Replace-Text "Yet, the intricate tapestry of quantum reality is painted with enigmatic paradoxes and perplexing dilemmas" "As we unravel the mysteries of the past, we are confronted by enigmatic figures who shaped the course of events, leaving an indelible mark on our world"

Replace-Text " The infamous Schrodinger's cat, living and dead at the same time, epitomizes the paradoxical nature of quantum superposition" " From kings and queens to pioneers and philosophers, their lives are entangled in a web of alliances, betrayals, and fortuitous coincidences"

Replace-Text " The act of observation, like a sorcerer's incantation, collapses the quantum waveform, forcing the cat into a definitive state of existence or nonexistence" " Their decisions, shaped by personal ambitions and societal pressures, ripple through time, influencing future generations in ways they could never have foreseen"

Replace-Text " But what is the catalyst that triggers this collapse? Does the conscious mind of an observer play an active role in shaping quantum outcomes? Or is it a more fundamental process, guided by the laws of probability?" " Yet, in the grand symphony of history, even the most towering individuals are but players, part of a larger ensemble where innumerable lives intersect and destinies intertwine."

Replace-Text "The pursuit of answers to these perplexing quandaries propels us on an intellectual odyssey that transcends the boundaries of discipline" "History is a canvas painted with countless brushstrokes, each contributing to the overall masterpiece"

Replace-Text " From the theoretical frameworks of quantum mechanics to the experimental frontiers of particle physics laboratories, scientists labor tirelessly to illuminate the hidden mechanisms that govern the quantum realm" " It is a kaleidoscope of cultures and beliefs, a symphony of voices recounting events both harrowing and sublime"

Replace-Text " Each new discovery, like a flash of lightning, illuminates another corner of this enigmatic universe, revealing tantalizing glimpses of a reality governed by laws that defy our everyday experiences" " It is in the study of this mosaic that we discover the profound interconnectedness of all human existence. For we are but the sum of our ancestors' experiences, shaped by their triumphs and missteps, failures and ingenuities. In understanding the past, we gain a keener appreciation for the present and a renewed sense of responsibility for the future"

# --- Summary paragraph ---
Replace-Text "Our exploration into the quantum realm unveils a tapestry of paradoxes and perplexities that defy conventional understanding" "History is a tapestry of human experience, woven from the threads of countless lives"

Replace-Text " Quantum entanglement and superposition challenge our notions of locality and causality" " It is a story of triumphs and tragedies, of heroism and folly, told through the voices of individuals and civilizations"

Replace-Text " The infamous Schrodinger's cat exemplifies the enigmatic nature of quantum superposition, prompting profound questions about the role of observers and consciousness" " In studying history, we gain a deeper understanding of ourselves and the world around us"

Replace-Text " The quest for answers to these riddles propels us on an intellectual journey, spanning disciplines and pushing the boundaries of human knowledge. The mysteries of the quantum realm beckon us, whispering promises of technological advancements and a deeper understanding of the very fabric of reality" " It challenges us to grapple with the complexities of the human condition and inspires us to strive for a better future"

# --- Add an empty paragraph at end of document (before sectPr) ---
# Appending a bare carriage return to the end of the document's text range
# splits off a brand-new, completely empty paragraph (no runs) -- matching
# `<w:p/>` -- instead of `InsertParagraphAfter()`, which leaves a stray
# empty run carrying the inherited run-formatting.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$lastPara.InsertAfter("`r")

Write-Output "done"
